$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.562.96'
$ws.Range('E2').Value = '  +7.44%  '
$ws.Range('D3').Value = '3.543.42'
$ws.Range('E3').Value = '  +9.84%  '
$ws.Range('D5').Value = "'193.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.12%  '
$ws.Range('D6').Value = "'558.35"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.26%  '
$ws.Range('D7').Value = '3.541.71'
$ws.Range('E7').Value = '  +9.84%  '
$ws.Range('D8').Value = "'0.611"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.43%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = "'0.643"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.62%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'0.152"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +17.03%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Value = "'56.74"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.22%  '
$ws.Range('D13').Value = "'0.0000273"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.02%  '
$ws.Range('D14').Value = "'9.54"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.61%  '
$ws.Range('D15').Value = '4.114.98'
$ws.Range('E15').Value = '  +9.84%  '
$ws.Range('D16').Value = '3.550.47'
$ws.Range('E16').Value = '  +10.04%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.122"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.02%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '67.691.65'
$ws.Range('E18').Value = '  +7.65%  '
$ws.Range('D19').Value = "'18.45"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.90%  '
$ws.Range('D20').Value = "'11.96"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.19%  '
$ws.Range('D21').Value = "'1.01"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.35%  '
$ws.Range('D22').Value = "'408.70"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +11.94%  '
$ws.Range('D23').Value = "'3.99"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.60%  '
$ws.Range('D24').Value = "'85.87"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.13%  '
$ws.Range('D25').Value = "'4.25"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.15%  '
$ws.Range('D26').Value = "'11.37"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.39%  '
$ws.Range('E27').Value = '  +14.54%  '
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').Value = "'12.09"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.12%  '
$ws.Range('D30').Value = "'8.85"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.94%  '
$ws.Range('D31').Value = "'30.64"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.92%  '
$ws.Range('D32').Value = "'689.15"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('D33').Value = "'6.85"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.05%  '
$ws.Range('D34').Value = "'11.86"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.33%  '
$ws.Range('E35').Value = '  +8.72%  '
$ws.Range('D36').Value = "'60.83"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.25%  '
$ws.Range('D37').Value = "'39.24"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.86%  '
$ws.Range('D38').Value = '0.0₃0828'
$ws.Range('E38').Value = '  +15.66%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  +7.35%  '
$ws.Range('E41').Value = '  +14.74%  '
$ws.Range('E42').Value = '  +21.18%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').Value = "'3.05"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +18.43%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'0.999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = "'2.69"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.17%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '3.037.25'
$ws.Range('E46').Value = '  +6.40%  '
$ws.Range('D47').Value = "'0.0422"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.89%  '
$ws.Range('D48').Value = "'3.25"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.62%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = "'2.74"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'9.11"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +19.34%  '
$ws.Range('E51').Value = '  +7.05%  '
